$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and G keep their existing text format so numeric-looking
# strings (prices, hour values) are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "239.64"
$ws.Range("G2").Value = "23"
$ws.Range("D3").Value = "21.71"
$ws.Range("G3").Value = "23"
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D4").Value = "5.389"
$ws.Range("E4").Value = "3HuobiTokenHT"
$ws.Range("G4").Value = "23"
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D5").Value = "0.05561"
$ws.Range("E5").Value = "4CronosCRO"
$ws.Range("G5").Value = "23"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "6.467"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("G6").Value = "23"
$ws.Range("D7").Value = "3.362"
$ws.Range("G7").Value = "23"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.8053"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("G8").Value = "23"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "1.072"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("G9").Value = "23"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1410"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").Value = "23"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07333"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Value = "23"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03286"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "23"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.02945"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Value = "23"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09247"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").Value = "23"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001643"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("G15").Value = "23"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.243"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("G16").Value = "23"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04774"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").Value = "23"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005698"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "23"
$ws.Range("D19").Value = "0.006246"
$ws.Range("G19").Value = "23"
$ws.Range("D20").Value = "0.001045"
$ws.Range("G20").Value = "23"
$ws.Range("D21").Value = "0.003795"
$ws.Range("G21").Value = "23"
$ws.Range("D22").Value = "0.0001496"
$ws.Range("G22").Value = "23"
$ws.Range("D23").Value = "0.0004172"
$ws.Range("G23").Value = "23"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "3.974"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("G24").Value = "23"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "2.203"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("G25").Value = "23"
$ws.Range("G26").Value = "23"
$ws.Range("G27").Value = "23"
$ws.Range("G28").Value = "23"
$ws.Range("G29").Value = "23"
$ws.Range("G30").Value = "23"
$ws.Range("G31").Value = "23"
$ws.Range("G32").Value = "23"
$ws.Range("G33").Value = "23"
$ws.Range("G34").Value = "23"
$ws.Range("G35").Value = "23"
$ws.Range("G36").Value = "23"
$ws.Range("G37").Value = "23"
$ws.Range("G38").Value = "23"
$ws.Range("G39").Value = "23"
$ws.Range("D40").Value = "0.04175"
$ws.Range("G40").Value = "23"
$ws.Range("D41").Value = "0.006980"
$ws.Range("G41").Value = "23"
$ws.Range("D42").Value = "0.003492"
$ws.Range("G42").Value = "23"
$ws.Range("D43").Value = "0.1042"
$ws.Range("G43").Value = "23"
$ws.Range("D44").Value = "0.009772"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("G44").Value = "23"
$ws.Range("D45").Value = "0.00005429"
$ws.Range("G45").Value = "23"
$ws.Range("D46").Value = "0.00000000748"
$ws.Range("G46").Value = "23"
$ws.Range("D47").Value = "0.6785"
$ws.Range("G47").Value = "23"
$ws.Range("D48").Value = "0.03111"
$ws.Range("G48").Value = "23"
$ws.Range("D49").Value = "0.00002095"
$ws.Range("G49").Value = "23"
$ws.Range("D50").Value = "0.01008"
$ws.Range("G50").Value = "23"
$ws.Range("G51").Value = "23"
